$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New bookmark rows to append (user_id, recipe_id, created_at)
# recipe_id values look numeric ("1", "14", "45") but must be stored as TEXT,
# matching the existing B column (e.g. B2 = "112140" stored as text).
$rows = @(
    @{ Row = 3; UserId = 4; RecipeId = "1";  Created = 45996.82861790509 },
    @{ Row = 4; UserId = 4; RecipeId = "14"; Created = 45996.84202701389 },
    @{ Row = 5; UserId = 5; RecipeId = "45"; Created = 45996.85839129629 },
    @{ Row = 6; UserId = 5; RecipeId = "14"; Created = 45996.858683252314 }
)

foreach ($r in $rows) {
    $row = $r.Row

    # Column A: plain numeric user_id
    $ws.Cells.Item($row, 1).Value = $r.UserId

    # Column B: recipe_id stored as text even though it looks like a number.
    # A leading apostrophe forces text entry (like typing '1 into Excel);
    # resetting the style back to Normal afterwards drops the quote-prefix
    # formatting so the cell ends up as a plain shared-string cell, just
    # like the other text cells in this sheet.
    $ws.Cells.Item($row, 2).Value = "'" + $r.RecipeId
    $ws.Cells.Item($row, 2).Style = "Normal"

    # Column C: created_at date/time value, using the same date number
    # format already used by C2 (style index 1 / numFmtId 14).
    $ws.Cells.Item($row, 3).Value = $r.Created
    $ws.Cells.Item($row, 3).NumberFormat = "m/d/yy"
}

Write-Output "Appended $($rows.Count) bookmark rows"
